$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.076.25"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "2.596.27"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'519.41"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'139.38"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "2.620.43"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").Value = "'6.48"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "3.061.97"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "59.057.80"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "'20.39"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "2.614.58"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "'339.70"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'4.31"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "'10.17"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'6.47"
$ws.Range("E22").Value = "  +5.90%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'66.31"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "'0.404"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "'7.04"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").Value = "'18.74"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'149.37"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "'36.33"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "'0.832"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.46"
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'0.835"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'276.56"
$ws.Range("E42").Value = "  +5.97%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "'10.73"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("D46").Value = "'0.0950"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'0.0520"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'18.53"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "1.987.09"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'4.59"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "'0.0220"
$ws.Range("E51").Value = "  -0.82%  "
